$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "THIS IS THE COOLEST SPREADSHEET EVER!" spelled out letter-by-letter
# across a block of cells (rows 7-10).
$letters = [ordered]@{
    "F7" = "T"; "G7" = "H"; "H7" = "I"; "I7" = "S"; "K7" = "I"; "L7" = "S";
    "D8" = "T"; "E8" = "H"; "F8" = "E"; "H8" = "C"; "I8" = "O"; "J8" = "O"; "K8" = "L"; "L8" = "E"; "M8" = "S"; "N8" = "T";
    "D9" = "S"; "E9" = "P"; "F9" = "R"; "G9" = "E"; "H9" = "A"; "I9" = "D"; "J9" = "S"; "K9" = "H"; "L9" = "E"; "M9" = "E"; "N9" = "T";
    "G10" = "E"; "H10" = "V"; "I10" = "E"; "J10" = "R"; "K10" = "!";
}

foreach ($addr in $letters.Keys) {
    $ws.Range($addr).Value = $letters[$addr]
}

$ws.Range("L10").Select()
